$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "general": add a new "fao_years" row, and add raw-data rows for
# genebank_file / upov_file / gbif_research_file / sgsv_file / sgsv_fields /
# primary_region_file / primary_region_fields (supporting "sgsv" and
# "primary region" data sources).
# ---------------------------------------------------------------------------
$wsGeneral = $wb.Worksheets.Item("general")

# New row 9: fao_years
$wsGeneral.Rows.Item(9).Insert()
$wsGeneral.Range("A9").Value = "fao_years"
$wsGeneral.Range("B9").Value = "'2015,2016,2017,2018"

# New row 17: genebank_file (genebank_fields already exists right after)
$wsGeneral.Rows.Item(17).Insert()
$wsGeneral.Range("A17").Value = "genebank_file"
$wsGeneral.Range("B17").Value = "genebank_collection.csv"

# New row 19: upov_file (upov_fields already exists right after)
$wsGeneral.Rows.Item(19).Insert()
$wsGeneral.Range("A19").Value = "upov_file"
$wsGeneral.Range("B19").Value = "upov_varietal_release.csv"

# New row 21: gbif_research_file (gbif_research_fields already exists right after)
$wsGeneral.Rows.Item(21).Insert()
$wsGeneral.Range("A21").Value = "gbif_research_file"
$wsGeneral.Range("B21").Value = "gbif_research_supply.csv"

# New rows 23-26: sgsv_file / sgsv_fields / primary_region_file / primary_region_fields
$wsGeneral.Range("A23:A26").EntireRow.Insert()
$wsGeneral.Range("A23").Value = "sgsv_file"
$wsGeneral.Range("B23").Value = "accessions_sgsv.csv"
$wsGeneral.Range("A24").Value = "sgsv_fields"
$wsGeneral.Range("B24").Value = "genus_accessions_sgsv,species_accessions_sgsv"
$wsGeneral.Range("A25").Value = "primary_region_file"
$wsGeneral.Range("B25").Value = "accessions_primaryregion.csv"
$wsGeneral.Range("A26").Value = "primary_region_fields"
$wsGeneral.Range("B26").Value = "genus_accessions_primaryregion,species_accessions_primaryregion"

# ---------------------------------------------------------------------------
# Sheet "downloads": cursor moved to A17.
# ---------------------------------------------------------------------------
$wsDownloads = $wb.Worksheets.Item("downloads")
$wsDownloads.Range("A17").Select()

# ---------------------------------------------------------------------------
# Sheet "indicator": fill in source/step/folder/file/element for the two
# "accessions_sgsv" rows (96-97), mirroring the other raw-data blocks.
# ---------------------------------------------------------------------------
$wsIndicator = $wb.Worksheets.Item("indicator")

$wsIndicator.Range("A96").Value = "sgsv"
$wsIndicator.Range("B96").Value = "'01"
$wsIndicator.Range("C96").Value = "SM"
$wsIndicator.Range("D96").Value = "accessions_sgsv.csv"
$wsIndicator.Range("E96").Value = "genus_accessions_sgsv"

$wsIndicator.Range("A97").Value = "sgsv"
$wsIndicator.Range("B97").Value = "'01"
$wsIndicator.Range("C97").Value = "SM"
$wsIndicator.Range("D97").Value = "accessions_sgsv.csv"
$wsIndicator.Range("E97").Value = "species_accessions_sgsv"

$wsIndicator.Range("G95").Select()

# ---------------------------------------------------------------------------
# Make "general" the active sheet/selection, as in the saved workbook.
# ---------------------------------------------------------------------------
$wsGeneral.Activate()
$wsGeneral.Range("A26").Select()
